# Update "想去人数" (interested-count) values in column F on the
# "展览" and "全部类型" sheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 1349
$ws1.Range("F16").Value = 8304
$ws1.Range("F23").Value = 571
$ws1.Range("F25").Value = 1148
$ws1.Range("F28").Value = 1684
$ws1.Range("F31").Value = 1913
$ws1.Range("F40").Value = 391

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 1349
$ws4.Range("F17").Value = 8304
$ws4.Range("F25").Value = 571
$ws4.Range("F27").Value = 1148
$ws4.Range("F30").Value = 1684
$ws4.Range("F32").Value = 1913
$ws4.Range("F41").Value = 391
